# RPA_RentOffersGenerator/results.xlsx
# Commit: "Broke Main workflow into more workflows - renamed activities"
#
# The underlying run produced 5 new rental-offer URLs (the RPA scrape
# refreshed the source listings), widened column C to fit the longest
# new URL, and left the selection/cursor parked at A11 (one row below
# a blank spacer row that was touched during the session).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Refresh the 5 scraped listing URLs in column C (C1:C5) ------------
$ws.Range("C1").Value = "https://www.imobiliare.ro/inchirieri-apartamente/cluj-napoca/buna-ziua/apartament-de-inchiriat-2-camere-XBB800047"
$ws.Range("C2").Value = "https://www.imobiliare.ro/inchirieri-apartamente/cluj-napoca/gheorgheni/apartament-de-inchiriat-2-camere-X8FG10042"
$ws.Range("C3").Value = "https://www.imobiliare.ro/inchirieri-apartamente/cluj-napoca/gheorgheni/apartament-de-inchiriat-2-camere-XBNT102IT"
$ws.Range("C4").Value = "https://www.imobiliare.ro/inchirieri-apartamente/cluj-napoca/europa/apartament-de-inchiriat-2-camere-X8FG1005O"
$ws.Range("C5").Value = "https://www.imobiliare.ro/inchirieri-apartamente/cluj-napoca/gheorgheni/apartament-de-inchiriat-2-camere-X8FG1005P"

# --- 2. Column C grew wider (new URLs are longer) --------------------------
$ws.Columns.Item(3).ColumnWidth = 110

# --- 3. Touch row 11 (a blank spacer row below the data) so it is present
#        in the sheet without any cell content ------------------------------
$ws.Rows.Item(11).Hidden = $true
$ws.Rows.Item(11).Hidden = $false

# --- 4. Leave the cursor / selection parked at A11 --------------------------
$ws.Range("A11").Select()
